{"js": "// Insert a new run containing \"dsaf\" immediately before the first run of\n// the document's first paragraph (which holds \"\u0432\u044b\u0430\u0430\u044b\u0432\"), giving the new\n// run the same run properties (rFonts/color/spacing/position/sz/shd) as\n// the existing run it precedes.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = body.paragraphs.items[0];\n\n// Build a FlatOPC-wrapped OOXML fragment for the new run so it lands as\n// its own <w:r> (rather than being merged character-by-character into the\n// neighbouring run's <w:t>).\nconst runOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r>' +\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Calibri\" w:eastAsia=\"Calibri\"/>' +\n  '<w:color w:val=\"auto\"/>' +\n  '<w:spacing w:val=\"0\"/>' +\n  '<w:position w:val=\"0\"/>' +\n  '<w:sz w:val=\"22\"/>' +\n  '<w:shd w:fill=\"auto\" w:val=\"clear\"/>' +\n  '</w:rPr>' +\n  '<w:t xml:space=\"preserve\">dsaf</w:t>' +\n  '</w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nparagraph.getRange(\"Start\").insertOoxml(runOoxml, \"Start\");\nawait context.sync();\n", "ps1": "# Insert a new run containing \"dsaf\" immediately before the first run of\n# the document's first paragraph (which holds \"\u0432\u044b\u0430\u0430\u044b\u0432\"), giving the new\n# run the same run properties (rFonts/color/spacing/position/sz/shd) as\n# the existing run it precedes.\n\n$d = $word.ActiveDocument\n$p = $d.Paragraphs.Item(1)\n\n$r = $p.Range\n$r.Collapse(1)  # wdCollapseStart -> caret at the very start of the paragraph\n\n# FlatOPC-wrapped OOXML fragment so the inserted text lands as its own\n# <w:r> (rather than being merged into the neighbouring run's <w:t>).\n$xml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Calibri\" w:eastAsia=\"Calibri\"/><w:color w:val=\"auto\"/><w:spacing w:val=\"0\"/><w:position w:val=\"0\"/><w:sz w:val=\"22\"/><w:shd w:fill=\"auto\" w:val=\"clear\"/></w:rPr><w:t xml:space=\"preserve\">dsaf</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n\n$r.InsertXML($xml)\n"}
